$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '37.110.83'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -0.11%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.047.46'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -1.29%  '

$ws.Range('E4').Value = '  +0.03%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '248.08'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.97%  '

$ws.Range('E6').Value = '  -2.02%  '

$ws.Range('E7').Value = '  -0.02%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '56.84'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -3.95%  '

$ws.Range('E9').Value = '  -2.38%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0773'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -3.18%  '

$ws.Range('E11').Value = '  +0.02%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '15.76'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -3.04%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.862'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +5.11%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '2.345.51'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -1.18%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '5.67'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +2.50%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '2.052.10'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -0.94%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '17.81'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +13.27%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '37.079.03'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +0.08%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '74.44'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -0.41%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.0₃0888'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -4.76%  '

$ws.Range('E21').Value = '  -2.29%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '236.69'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -1.52%  '

$ws.Range('E23').Value = '  -0.03%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.46'
$ws.Range('D24').Style = "Normal"

$ws.Range('E25').Value = '  -4.71%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.43'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +0.68%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '168.93'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -0.43%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '19.99'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -1.62%  '

$ws.Range('E29').Value = '  -1.93%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '4.82'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +0.03%  '

$ws.Range('E31').Value = '  -1.94%  '

$ws.Range('E32').Value = '  -3.17%  '

$ws.Range('E33').Value = '  -0.63%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.0886'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -2.21%  '

$ws.Range('E35').Value = '  +0.02%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.24'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -2.79%  '

$ws.Range('E37').Value = '  +0.15%  '

$ws.Range('E38').Value = '  -2.11%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '3.16'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +12.05%  '

$ws.Range('E40').Value = '  +16.48%  '

$ws.Range('B41').Value = 'Cronos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0974'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -17.38%  '

$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.0222'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -2.39%  '

$ws.Range('B43').Value = 'InjectiveProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '17.13'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -3.42%  '

$ws.Range('B44').Value = 'ARBITRUM'
$ws.Range('C44').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.13'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -2.49%  '

$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '95.46'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -3.55%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.42'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -2.17%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.265.05'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -2.97%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.86'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -2.80%  '

$ws.Range('E49').Value = '  -2.27%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.230.35'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -0.96%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '43.57'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -1.24%  '
